# Update "想去人数" (F) and "最低票价" (G) figures on the 展览 and 全部类型
# sheets to match the refreshed bilibili export (gh-pages output regenerated
# at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value = 182      # F3: 181 -> 182
$ws1.Cells.Item(3, 7).Value = 39.9     # G3: 55  -> 39.9
$ws1.Cells.Item(4, 6).Value = 178      # F4: 177 -> 178
$ws1.Cells.Item(5, 6).Value = 5036     # F5: 5025 -> 5036
$ws1.Cells.Item(7, 6).Value = 46       # F7: 45 -> 46
$ws1.Cells.Item(9, 6).Value = 557      # F9: 555 -> 557
$ws1.Cells.Item(10, 6).Value = 516     # F10: 515 -> 516
$ws1.Cells.Item(11, 6).Value = 1039    # F11: 1038 -> 1039
$ws1.Cells.Item(13, 6).Value = 1404    # F13: 1402 -> 1404
$ws1.Cells.Item(14, 6).Value = 3703    # F14: 3687 -> 3703
$ws1.Cells.Item(16, 6).Value = 140     # F16: 137 -> 140
$ws1.Cells.Item(18, 6).Value = 84      # F18: 83 -> 84
$ws1.Cells.Item(19, 6).Value = 2701    # F19: 2694 -> 2701
$ws1.Cells.Item(21, 6).Value = 22      # F21: 19 -> 22
$ws1.Cells.Item(22, 6).Value = 90      # F22: 89 -> 90
$ws1.Cells.Item(24, 6).Value = 180     # F24: 178 -> 180
$ws1.Cells.Item(25, 6).Value = 61      # F25: 59 -> 61

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 6).Value = 182      # F3: 181 -> 182
$ws4.Cells.Item(3, 7).Value = 39.9     # G3: 55  -> 39.9
$ws4.Cells.Item(4, 6).Value = 178      # F4: 177 -> 178
$ws4.Cells.Item(6, 6).Value = 5036     # F6: 5025 -> 5036
$ws4.Cells.Item(8, 6).Value = 46       # F8: 45 -> 46
$ws4.Cells.Item(10, 6).Value = 557     # F10: 555 -> 557
$ws4.Cells.Item(11, 6).Value = 516     # F11: 515 -> 516
$ws4.Cells.Item(12, 6).Value = 1039    # F12: 1038 -> 1039
$ws4.Cells.Item(14, 6).Value = 1404    # F14: 1402 -> 1404
$ws4.Cells.Item(15, 6).Value = 3703    # F15: 3687 -> 3703
$ws4.Cells.Item(17, 6).Value = 140     # F17: 137 -> 140
$ws4.Cells.Item(19, 6).Value = 84      # F19: 83 -> 84
$ws4.Cells.Item(20, 6).Value = 2701    # F20: 2694 -> 2701
$ws4.Cells.Item(22, 6).Value = 22      # F22: 19 -> 22
$ws4.Cells.Item(23, 6).Value = 90      # F23: 89 -> 90
$ws4.Cells.Item(25, 6).Value = 180     # F25: 178 -> 180
$ws4.Cells.Item(26, 6).Value = 61      # F26: 59 -> 61
